$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings are not auto-converted to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.055.85'
$ws.Range("E2").Value = '  +6.27%  '
$ws.Range("D3").Value = '2.734.42'
$ws.Range("E3").Value = '  +4.64%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '591.51'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").Value = '153.30'
$ws.Range("E6").Value = '  +7.02%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '0.611'
$ws.Range("E8").Value = '  +2.25%  '
$ws.Range("D9").Value = '2.767.22'
$ws.Range("E9").Value = '  +5.53%  '
$ws.Range("D10").Value = '6.72'
$ws.Range("E10").Value = '  +3.35%  '
$ws.Range("E11").Value = '  +7.14%  '
$ws.Range("E12").Value = '  +3.67%  '
$ws.Range("E13").Value = '  +2.05%  '
$ws.Range("D14").Value = '3.224.26'
$ws.Range("E15").Value = '  +5.93%  '
$ws.Range("D16").Value = '63.901.95'
$ws.Range("E16").Value = '  +6.02%  '
$ws.Range("D17").Value = '0.0000152'
$ws.Range("E17").Value = '  +8.61%  '
$ws.Range("D18").Value = '2.762.28'
$ws.Range("E18").Value = '  +5.55%  '
$ws.Range("D19").Value = '12.08'
$ws.Range("E19").Value = '  +5.38%  '
$ws.Range("D20").Value = '4.90'
$ws.Range("E20").Value = '  +4.49%  '
$ws.Range("D21").Value = '365.16'
$ws.Range("E21").Value = '  +5.16%  '
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").Value = '0.538'
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D25").Value = '66.22'
$ws.Range("E25").Value = '  +3.99%  '
$ws.Range("E26").Value = '  +5.14%  '
$ws.Range("D27").Value = '8.67'
$ws.Range("E27").Value = '  +7.86%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '0.0₃0908'
$ws.Range("E29").Value = '  +13.55%  '
$ws.Range("E30").Value = '  +4.44%  '
$ws.Range("E31").Value = '  +9.21%  '
$ws.Range("D32").Value = '172.84'
$ws.Range("E32").Value = '  +2.53%  '
$ws.Range("D33").Value = '1.20'
$ws.Range("E33").Value = '  +18.33%  '
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '20.62'
$ws.Range("E35").Value = '  +5.69%  '
$ws.Range("D36").Value = '4.85'
$ws.Range("E36").Value = '  +13.12%  '
$ws.Range("E37").Value = '  +9.88%  '
$ws.Range("D38").Value = '1.78'
$ws.Range("E38").Value = '  +9.02%  '
$ws.Range("D39").Value = '1.02'
$ws.Range("E39").Value = '  +19.41%  '
$ws.Range("D40").Value = '348.63'
$ws.Range("E40").Value = '  +9.02%  '
$ws.Range("E41").Value = '  +7.11%  '
$ws.Range("D42").Value = '38.97'
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '22.11'
$ws.Range("E43").Value = '  +10.36%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '5.60'
$ws.Range("E44").Value = '  +11.36%  '
$ws.Range("D45").Value = '143.72'
$ws.Range("E45").Value = '  +5.84%  '
$ws.Range("D46").Value = '22.18'
$ws.Range("E46").Value = '  +10.74%  '
$ws.Range("E47").Value = '  +7.02%  '
$ws.Range("D48").Value = '0.649'
$ws.Range("E48").Value = '  +6.60%  '
$ws.Range("D49").Value = '0.0258'
$ws.Range("E49").Value = '  +7.14%  '
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").Value = '2.176.33'

# Restore default style on column D (remove the temporary text-format override)
$ws.Range("D2:D51").Style = "Normal"
